$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: the workday's end time moved earlier (16 -> 13) and a description
# of the day's work was added. Time worked (D15) recalculates automatically
# via the existing shared formula (=C15-B15).
$ws.Range("C15").Value = 13
$ws.Range("E15").Value = "Testing new token fetching and updating, planning out what and how to display from character info on frontend"

# Row 16: start/end time filled in (previously blank) and a description
# added. D16 (=C16-B16) recalculates automatically.
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 15
$ws.Range("E16").Value = "Adding resource redirects to each respective class, styling for ListViews when character list is implemented"

# The active selection moved from E14 to E16.
$ws.Range("E16").Select() | Out-Null
